# Update points 09876543 -> 0.00
#
# Original row 54 held phone "09876543" (text, leading zero) with 120 points.
# This edit:
#   1. Converts the existing row 54's phone to a plain number (9876543),
#      keeping its 120 points as-is.
#   2. Appends a new row 55 for phone "09876543" (text, leading zero kept)
#      with 0 points and a blank birthday, mirroring the neighbouring rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 54: "09876543" (text) -> 9876543 (number); points (C54) unchanged ---
$ws.Cells.Item(54, 1).Value = 9876543

# --- Row 55 (new): phone stays text w/ leading zero, points reset to 0 ---
# Leading apostrophe forces text interpretation so the leading zero survives.
$ws.Cells.Item(55, 1).Value = "'09876543"
$ws.Cells.Item(55, 2).Value = ""
$ws.Cells.Item(55, 3).Value = 0
